$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce / Rubber
$ws.Range("H11").Value = 1102.8889
$ws.Range("I11").Value = 1102.8889
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1102.8889
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -962.8888999999999

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 2316.5
$ws.Range("I40").Value = 2828.5
$ws.Range("J40").Value = 1676.5
$ws.Range("K40").Value = 2828.5
$ws.Range("L40").Value = 1676.5
$ws.Range("M40").Value = -2653.5
$ws.Range("N40").Value = -2026.5

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1367.1277
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 1410.1111
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 4230.3333
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -6446.3333

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1377.8125
$ws.Range("I132").Value = 1187.3684
$ws.Range("J132").Value = 2928.5715
$ws.Range("K132").Value = 3562.1052
$ws.Range("L132").Value = 8785.7145
$ws.Range("M132").Value = -1032.1052
$ws.Range("N132").Value = -13845.7145

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2703.56
$ws.Range("I137").Value = 1469.3889
$ws.Range("J137").Value = 5877.143
$ws.Range("K137").Value = 4408.1667
$ws.Range("L137").Value = 17631.429
$ws.Range("M137").Value = -1858.1667
$ws.Range("N137").Value = -22731.429

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3790.3684
$ws.Range("I138").Value = 2442.5
$ws.Range("J138").Value = 4271.75
$ws.Range("K138").Value = 7327.5
$ws.Range("L138").Value = 12815.25
$ws.Range("M138").Value = -2187.5
$ws.Range("N138").Value = -23095.25

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 63831.58
$ws.Range("I74").Value = 71265.664
$ws.Range("J74").Value = 16439.25
$ws.Range("K74").Value = 71265.664
$ws.Range("L74").Value = 16439.25
$ws.Range("M74").Value = -70391.664
$ws.Range("N74").Value = -18187.25

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 63831.58
$ws.Range("I77").Value = 71265.664
$ws.Range("J77").Value = 16439.25
$ws.Range("K77").Value = 356328.32
$ws.Range("L77").Value = 82196.25
$ws.Range("M77").Value = -351960.32
$ws.Range("N77").Value = -90932.25

# Row 109: A Head of Demand / Deepgold Helm of Fending
$ws.Range("H109").Value = 34626.2
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 34626.2
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 34626.2
$ws.Range("N109").Value = -37400.2

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 28195.666
$ws.Range("I134").Value = 2531.1143
$ws.Range("J134").Value = 252760.5
$ws.Range("K134").Value = 7593.342900000001
$ws.Range("L134").Value = 758281.5
$ws.Range("M134").Value = -5058.342900000001
$ws.Range("N134").Value = -763351.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1997.5
$ws.Range("I16").Value = 1996.6666
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1996.6666
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2574

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2299.5881
$ws.Range("I31").Value = 1626.509
$ws.Range("J31").Value = 5147.231
$ws.Range("K31").Value = 1626.509
$ws.Range("L31").Value = 5147.231
$ws.Range("M31").Value = -1331.509
$ws.Range("N31").Value = -5737.231

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2299.5881
$ws.Range("I34").Value = 1626.509
$ws.Range("J34").Value = 5147.231
$ws.Range("K34").Value = 1626.509
$ws.Range("L34").Value = 5147.231
$ws.Range("M34").Value = -1424.509
$ws.Range("N34").Value = -5551.231

# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 1096.7727
$ws.Range("I94").Value = 1232.75
$ws.Range("J94").Value = 1019.0714
$ws.Range("K94").Value = 1232.75
$ws.Range("L94").Value = 1019.0714
$ws.Range("M94").Value = -781.75
$ws.Range("N94").Value = -1921.0714

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 426.1
$ws.Range("I105").Value = 450
$ws.Range("J105").Value = 211
$ws.Range("K105").Value = 450
$ws.Range("L105").Value = 211
$ws.Range("N105").Value = -3705

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1997.5
$ws.Range("I113").Value = 1996.6666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1996.6666
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340

# Row 124: Earring Awakening / Palm Ear Cuffs of Fending
$ws.Range("H124").Value = 53000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 53000
$ws.Range("K124").Value = 0
$ws.Range("N124").Value = -57910

$ws = $wb.Worksheets.Item("CUL")
# Row 63: The Next to Last Supper / Stuffed Cabbage Rolls
$ws.Range("H63").Value = 3462.077
$ws.Range("I63").Value = 3163
$ws.Range("J63").Value = 3551.8
$ws.Range("K63").Value = 9489
$ws.Range("L63").Value = 10655.4
$ws.Range("M63").Value = -8740
$ws.Range("N63").Value = -12153.4

# Row 66: Nostalgia through the Stomach (L) / Stuffed Cabbage Rolls
$ws.Range("H66").Value = 3462.077
$ws.Range("I66").Value = 3163
$ws.Range("J66").Value = 3551.8
$ws.Range("K66").Value = 28467
$ws.Range("L66").Value = 31966.2
$ws.Range("M66").Value = -24723
$ws.Range("N66").Value = -39454.2

# Row 97: The Frier Never Lies / Cottonseed Oil
$ws.Range("H97").Value = 7133.9375
$ws.Range("I97").Value = 763.3333
$ws.Range("J97").Value = 10956.3
$ws.Range("K97").Value = 2289.9999
$ws.Range("L97").Value = 32868.89999999999
$ws.Range("M97").Value = -1793.9999
$ws.Range("N97").Value = -33860.89999999999

# Row 114: One Last Meal / Mushroom Saute
$ws.Range("H114").Value = 567.93335
$ws.Range("I114").Value = 460.83334
$ws.Range("J114").Value = 996.3333
$ws.Range("K114").Value = 1382.50002
$ws.Range("L114").Value = 2988.9999
$ws.Range("M114").Value = 1871.49998
$ws.Range("N114").Value = -9496.999899999999

# Row 117: A Good Omen / Peppered Popotoes
$ws.Range("H117").Value = 564.5
$ws.Range("I117").Value = 564.5
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1693.5
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1748.5

# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 1900.6333
$ws.Range("I129").Value = 2058
$ws.Range("J129").Value = 1780.2941
$ws.Range("K129").Value = 6174
$ws.Range("L129").Value = 5340.8823
$ws.Range("M129").Value = -1174
$ws.Range("N129").Value = -15340.8823

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 4199
$ws.Range("I107").Value = 398
$ws.Range("J107").Value = 8000
$ws.Range("K107").Value = 398
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = 1522
$ws.Range("N107").Value = -11840

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 16804.5
$ws.Range("I122").Value = 50000
$ws.Range("J122").Value = 5739.3335
$ws.Range("K122").Value = 150000
$ws.Range("L122").Value = 17218.0005
$ws.Range("M122").Value = -147550
$ws.Range("N122").Value = -22118.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3387.2
$ws.Range("I132").Value = 3191.0476
$ws.Range("J132").Value = 3681.4285
$ws.Range("K132").Value = 9573.1428
$ws.Range("L132").Value = 11044.2855
$ws.Range("M132").Value = -7043.1428

$ws = $wb.Worksheets.Item("WVR")
# Row 52: Party Animals / Linen Deerstalker
$ws.Range("H52").Value = 14700
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 19400
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 19400
$ws.Range("M52").Value = -9774
$ws.Range("N52").Value = -19852

# Row 58: Seeing It Through to the End / Woolen Smock
$ws.Range("H58").Value = 17950
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 17950
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 17950
$ws.Range("N58").Value = -18566
